# Straatroof.xlsx — "Testing for pull request"
#
# The Notification table's Message column gets a few entries reworded to use
# bracketed placeholders (to be filled in at runtime, e.g. "[geslacht]",
# "[bovenstuk]", "[onderstuk]"), one row's Postable flag flips from FALSE to
# TRUE, and the current selection moves to D1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Reworded notification messages (column D) — placeholders added for
# gender/clothing substitution.
$ws.Range("D3").Value  = "Een [geslacht] stal een tas van een vrouw."
$ws.Range("D7").Value  = "De[geslacht] rende hier gauw weg!"
$ws.Range("D10").Value = "Ik zag iemand hier wegrennen met een tas en een [bovenstuk] bovenstuk!"
$ws.Range("D14").Value = "De persoon die de tas dumpte had een [bovenstuk] shirt aan."
$ws.Range("D20").Value = "De overvaller loopt hier volgens mij. De [geslacht] draagt een [onderstuk] onderstuk! "

# Row 4 (Id=2, "Waar is dit gebeurd?") is now postable.
$ws.Range("G4").Value = $true

# Selection moves from C19 to D1.
$ws.Range("D1").Select()
